$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.427210000000001
$ws.Range("H2").Value = 28.28163
$ws.Range("I2").Value = 0.2188083857550241
$ws.Range("J2").Value = 0.2188083857550241
$ws.Range("M2").Value = 0.5134303333333333
$ws.Range("N2").Value = 1.540291
$ws.Range("O2").Value = 0.03326489761800302
$ws.Range("P2").Value = 0.03326489761800301
$ws.Range("Q2").Value = 4.840215572703333
$ws.Range("R2").Value = 43.56194015433
$ws.Range("S2").Value = 0.007278638550101388
$ws.Range("T2").Value = 0.007278638550101386
$ws.Range("G3").Value = 9.427210000000001
$ws.Range("H3").Value = 28.28163
$ws.Range("I3").Value = 0.2188083857550241
$ws.Range("J3").Value = 0.2188083857550241
$ws.Range("M3").Value = 1.626140333333333
$ws.Range("N3").Value = 4.878420999999999
$ws.Range("O3").Value = 0.1053568287437347
$ws.Range("P3").Value = 0.1053568287437347
$ws.Range("Q3").Value = 15.32996641180333
$ws.Range("R3").Value = 137.96969770623
$ws.Range("S3").Value = 0.02305295762568512
$ws.Range("T3").Value = 0.02305295762568512
$ws.Range("G4").Value = 9.427210000000001
$ws.Range("H4").Value = 28.28163
$ws.Range("I4").Value = 0.2188083857550241
$ws.Range("J4").Value = 0.2188083857550241
$ws.Range("M4").Value = 10.254745
$ws.Range("N4").Value = 30.764235
$ws.Range("O4").Value = 0.6643998618255804
$ws.Range("P4").Value = 0.6643998618255803
$ws.Range("Q4").Value = 96.67363461145
$ws.Range("R4").Value = 870.06271150305
$ws.Range("S4").Value = 0.1453762612619163
$ws.Range("T4").Value = 0.1453762612619163
$ws.Range("G5").Value = 9.427210000000001
$ws.Range("H5").Value = 28.28163
$ws.Range("I5").Value = 0.2188083857550241
$ws.Range("J5").Value = 0.2188083857550241
$ws.Range("M5").Value = 3.040282666666667
$ws.Range("N5").Value = 9.120848000000001
$ws.Range("O5").Value = 0.1969784118126819
$ws.Range("P5").Value = 0.1969784118126819
$ws.Range("Q5").Value = 28.66138315802667
$ws.Range("R5").Value = 257.95244842224
$ws.Range("S5").Value = 0.0431005283173213
$ws.Range("T5").Value = 0.0431005283173213
$ws.Range("I6").Value = 0.3808887290954196
$ws.Range("J6").Value = 0.3808887290954196
$ws.Range("M6").Value = 0.5134303333333333
$ws.Range("N6").Value = 1.540291
$ws.Range("O6").Value = 0.03326489761800302
$ws.Range("P6").Value = 0.03326489761800301
$ws.Range("Q6").Value = 8.425561715440333
$ws.Range("R6").Value = 75.83005543896299
$ws.Range("S6").Value = 0.01267022457721042
$ws.Range("T6").Value = 0.01267022457721042
$ws.Range("I7").Value = 0.3808887290954196
$ws.Range("J7").Value = 0.3808887290954196
$ws.Range("M7").Value = 1.626140333333333
$ws.Range("N7").Value = 4.878420999999999
$ws.Range("O7").Value = 0.1053568287437347
$ws.Range("P7").Value = 0.1053568287437347
$ws.Range("Q7").Value = 26.68550112245033
$ws.Range("S7").Value = 0.04012922860172489
$ws.Range("T7").Value = 0.04012922860172489
$ws.Range("I8").Value = 0.3808887290954196
$ws.Range("J8").Value = 0.3808887290954196
$ws.Range("M8").Value = 10.254745
$ws.Range("N8").Value = 30.764235
$ws.Range("O8").Value = 0.6643998618255804
$ws.Range("P8").Value = 0.6643998618255803
$ws.Range("Q8").Value = 168.283759770595
$ws.Range("R8").Value = 1514.553837935355
$ws.Range("S8").Value = 0.2530624189819177
$ws.Range("T8").Value = 0.2530624189819177
$ws.Range("I9").Value = 0.3808887290954196
$ws.Range("J9").Value = 0.3808887290954196
$ws.Range("M9").Value = 3.040282666666667
$ws.Range("N9").Value = 9.120848000000001
$ws.Range("O9").Value = 0.1969784118126819
$ws.Range("P9").Value = 0.1969784118126819
$ws.Range("Q9").Value = 49.89204489356266
$ws.Range("R9").Value = 449.028404042064
$ws.Range("S9").Value = 0.07502685693456659
$ws.Range("T9").Value = 0.07502685693456659
$ws.Range("G10").Value = 7.213061
$ws.Range("H10").Value = 21.639183
$ws.Range("I10").Value = 0.1674173200514808
$ws.Range("J10").Value = 0.1674173200514808
$ws.Range("M10").Value = 0.5134303333333333
$ws.Range("N10").Value = 1.540291
$ws.Range("O10").Value = 0.03326489761800302
$ws.Range("P10").Value = 0.03326489761800301
$ws.Range("Q10").Value = 3.703404313583667
$ws.Range("R10").Value = 33.33063882225299
$ws.Range("S10").Value = 0.005569120010992952
$ws.Range("T10").Value = 0.005569120010992951
$ws.Range("G11").Value = 7.213061
$ws.Range("H11").Value = 21.639183
$ws.Range("I11").Value = 0.1674173200514808
$ws.Range("J11").Value = 0.1674173200514808
$ws.Range("M11").Value = 1.626140333333333
$ws.Range("N11").Value = 4.878420999999999
$ws.Range("O11").Value = 0.1053568287437347
$ws.Range("P11").Value = 0.1053568287437347
$ws.Range("Q11").Value = 11.72944941889366
$ws.Range("R11").Value = 105.565044770043
$ws.Range("S11").Value = 0.01763855791739888
$ws.Range("T11").Value = 0.01763855791739888
$ws.Range("G12").Value = 7.213061
$ws.Range("H12").Value = 21.639183
$ws.Range("I12").Value = 0.1674173200514808
$ws.Range("J12").Value = 0.1674173200514808
$ws.Range("M12").Value = 10.254745
$ws.Range("N12").Value = 30.764235
$ws.Range("O12").Value = 0.6643998618255804
$ws.Range("P12").Value = 0.6643998618255803
$ws.Range("Q12").Value = 73.96810122444499
$ws.Range("R12").Value = 665.712911020005
$ws.Range("S12").Value = 0.1112320443094128
$ws.Range("T12").Value = 0.1112320443094128
$ws.Range("G13").Value = 7.213061
$ws.Range("H13").Value = 21.639183
$ws.Range("I13").Value = 0.1674173200514808
$ws.Range("J13").Value = 0.1674173200514808
$ws.Range("M13").Value = 3.040282666666667
$ws.Range("N13").Value = 9.120848000000001
$ws.Range("O13").Value = 0.1969784118126819
$ws.Range("P13").Value = 0.1969784118126819
$ws.Range("Q13").Value = 21.92974433190933
$ws.Range("R13").Value = 197.367698987184
$ws.Range("S13").Value = 0.03297759781367614
$ws.Range("T13").Value = 0.03297759781367614
$ws.Range("G14").Value = 10.03371566666667
$ws.Range("H14").Value = 30.101147
$ws.Range("I14").Value = 0.2328855650980756
$ws.Range("J14").Value = 0.2328855650980756
$ws.Range("M14").Value = 0.5134303333333333
$ws.Range("N14").Value = 1.540291
$ws.Range("O14").Value = 0.03326489761800302
$ws.Range("P14").Value = 0.03326489761800301
$ws.Range("Q14").Value = 5.151613979308556
$ws.Range("R14").Value = 46.364525813777
$ws.Range("S14").Value = 0.007746914479698262
$ws.Range("T14").Value = 0.00774691447969826
$ws.Range("G15").Value = 10.03371566666667
$ws.Range("H15").Value = 30.101147
$ws.Range("I15").Value = 0.2328855650980756
$ws.Range("J15").Value = 0.2328855650980756
$ws.Range("M15").Value = 1.626140333333333
$ws.Range("N15").Value = 4.878420999999999
$ws.Range("O15").Value = 0.1053568287437347
$ws.Range("P15").Value = 0.1053568287437347
$ws.Range("Q15").Value = 16.31622973876522
$ws.Range("R15").Value = 146.846067648887
$ws.Range("S15").Value = 0.02453608459892583
$ws.Range("T15").Value = 0.02453608459892583
$ws.Range("G16").Value = 10.03371566666667
$ws.Range("H16").Value = 30.101147
$ws.Range("I16").Value = 0.2328855650980756
$ws.Range("J16").Value = 0.2328855650980756
$ws.Range("M16").Value = 10.254745
$ws.Range("N16").Value = 30.764235
$ws.Range("O16").Value = 0.6643998618255804
$ws.Range("P16").Value = 0.6643998618255803
$ws.Range("Q16").Value = 102.8931955641717
$ws.Range("R16").Value = 926.0387600775451
$ws.Range("S16").Value = 0.1547291372723336
$ws.Range("T16").Value = 0.1547291372723336
$ws.Range("G17").Value = 10.03371566666667
$ws.Range("H17").Value = 30.101147
$ws.Range("I17").Value = 0.2328855650980756
$ws.Range("J17").Value = 0.2328855650980756
$ws.Range("M17").Value = 3.040282666666667
$ws.Range("N17").Value = 9.120848000000001
$ws.Range("O17").Value = 0.1969784118126819
$ws.Range("P17").Value = 0.1969784118126819
$ws.Range("Q17").Value = 30.50533182362845
$ws.Range("R17").Value = 274.547986412656
$ws.Range("S17").Value = 0.04587342874711787
$ws.Range("T17").Value = 0.04587342874711787
